$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date for 04cdfa6a-...md (shared between Overview and de-de sheets)
$wsOverview.Range("G3").Value = "2016-08-14 01:01:49"
$wsDeDe.Range("H3").Value = "2016-08-14 01:01:49"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 3
$wsZhCn.Range("H3").Value = "2016-08-14 01:01:41"
$wsZhCn.Range("K3").Value = "2016-08-14 01:02:15"

# de-de sheet: Correspond Handback DateTime for row 3
$wsDeDe.Range("K3").Value = "2016-08-14 01:02:26"
